# Weekly update: insert the newest week's "Cebollín" (Vega Monumental Concepción)
# price rows at the top of the data block (rows 12-13), pushing the older
# historical rows (old 12-36) down by two rows (new 14-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 12, shifting rows 12:36 down to 14:38.
$ws.Range("A12:R13").EntireRow.Insert()

# Populate the new row 12 ("Primera" quality) with this week's figures.
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 44679
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 100112037
$ws.Cells.Item(12, 7).Value = "Cebollín"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 600
$ws.Cells.Item(12, 12).Value = 700
$ws.Cells.Item(12, 13).Value = 650
$ws.Cells.Item(12, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(12, 15).Value = "Región de Ñuble"
$ws.Cells.Item(12, 16).Value = 108
$ws.Cells.Item(12, 17).Value = 6
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# Populate the new row 13 ("Segunda" quality) with this week's figures.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44679
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112037
$ws.Cells.Item(13, 7).Value = "Cebollín"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 500
$ws.Cells.Item(13, 12).Value = 500
$ws.Cells.Item(13, 13).Value = 500
$ws.Cells.Item(13, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(13, 15).Value = "Región de Ñuble"
$ws.Cells.Item(13, 16).Value = 83
$ws.Cells.Item(13, 17).Value = 6
$ws.Cells.Item(13, 18).Value = "Hortaliza"
